# HeadingDiscrimination_test.xlsx edit
# - Rename "increment" header (col I) to "step_plus"
# - Add new "step_mult" header in col J
# - Convert the "type" column (C) from a numeric code to a descriptive text:
#     1 (and the 0 on row 6) -> "constant"
#     2                      -> "varied"
# - Row 6: editable (D6) flips from 1 to 0
# - Column C width grows or shrink to fit the new text, bestFit turned off
# - Selection moves from E38 to E32
# - Used range grows from A1:I27 to A1:J27

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
$ws.Range("I1").Value = "step_plus"
$ws.Range("J1").Value = "step_mult"

# --- Column C: numeric type code -> text status ------------------------
# Rows 2-18, 25-26 hold a 0/1 (constant) or 2 (varied) code today;
# rows 19-24 and 27 stay on their -1 numeric sentinel and are untouched.
$constantRows = 2,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,25,26
foreach ($r in $constantRows) {
    $ws.Range("C$r").Value = "constant"
}
$ws.Range("C3").Value = "varied"

# Row 6 also had its "editable" flag (col D) flipped off
$ws.Range("D6").Value = 0

# --- Column C width (no longer auto-best-fit, now a bit wider) ---------
$ws.Columns.Item(3).ColumnWidth = 8.42

# --- Selection moves ----------------------------------------------------
$null = $ws.Range("E32").Select()
